$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 490
$ws.Range("G3").Value = "已售罄"
$ws.Range("F4").Value = 1462
$ws.Range("F5").Value = 748
$ws.Range("F6").Value = 170
$ws.Range("F7").Value = 39
$ws.Range("F8").Value = 1057
$ws.Range("F9").Value = 628
$ws.Range("F10").Value = 727
$ws.Range("F11").Value = 1256
$ws.Range("F12").Value = 250
$ws.Range("F13").Value = 983
$ws.Range("F14").Value = 41
$ws.Range("F15").Value = 171
$ws.Range("F16").Value = 29
$ws.Range("F17").Value = 367
$ws.Range("F20").Value = 509
$ws.Range("F21").Value = 530
$ws.Range("F22").Value = 718
$ws.Range("F23").Value = 206
$ws.Range("F24").Value = 143
$ws.Range("F25").Value = 349

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = "不可售"
$ws.Range("F4").Value = 950
$ws.Range("F6").Value = 180
$ws.Range("F7").Value = 12
$ws.Range("F9").Value = 55
$ws.Range("F10").Value = 552
$ws.Range("G10").Value = 580

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 154

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 490
$ws.Range("G4").Value = "已售罄"
$ws.Range("F5").Value = 154
$ws.Range("F6").Value = 1462
$ws.Range("F8").Value = 748
$ws.Range("F9").Value = 170
$ws.Range("F10").Value = 950
$ws.Range("F11").Value = 39
$ws.Range("F12").Value = 1057
$ws.Range("F13").Value = 628
$ws.Range("F14").Value = 727
$ws.Range("F15").Value = 1256
$ws.Range("F16").Value = 250
$ws.Range("F17").Value = 983
$ws.Range("F18").Value = 41
$ws.Range("F19").Value = 171
$ws.Range("F20").Value = 29
$ws.Range("F21").Value = 367
$ws.Range("F23").Value = 180
$ws.Range("F26").Value = 12
$ws.Range("F29").Value = 509
$ws.Range("F30").Value = 530
$ws.Range("F31").Value = 718
$ws.Range("F32").Value = 206
$ws.Range("F33").Value = 55
$ws.Range("F34").Value = 143
$ws.Range("F35").Value = 552
$ws.Range("G35").Value = 580
$ws.Range("F37").Value = 349
